# Update "想去人数" (F column) counts on three worksheets:
#   展览 (Exhibition)   - Worksheets.Item(1)
#   演出 (Show)         - Worksheets.Item(2)
#   全部类型 (All types) - Worksheets.Item(4)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3678
$ws1.Range("F6").Value = 437
$ws1.Range("F7").Value = 6
$ws1.Range("F8").Value = 5
$ws1.Range("F9").Value = 176
$ws1.Range("F11").Value = 77
$ws1.Range("F12").Value = 1362
$ws1.Range("F14").Value = 2076
$ws1.Range("F15").Value = 149

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 4

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3678
$ws4.Range("F6").Value = 437
$ws4.Range("F7").Value = 6
$ws4.Range("F8").Value = 5
$ws4.Range("F10").Value = 176
$ws4.Range("F12").Value = 77
$ws4.Range("F14").Value = 4
$ws4.Range("F15").Value = 1362
$ws4.Range("F17").Value = 2076
$ws4.Range("F18").Value = 149
